# Auto-update TW Market Data: 2026-02-20 07:50:50 UTC
#
# Sheet "1_Daily_Signals" (sheet1):
#   - rows 30/31 swap (2454 MediaTek <-> 2002 China Steel)
#   - rows 36/37 swap (2615 Wan Hai <-> 3260 ADATA)
#   - new row inserted at 47 for 2303 / UMC, pushing old rows 47-51 down to 48-52
#
# Sheet "2_21Day_Trend" (sheet2):
#   - new row inserted at 10 for 2303 / UMC, pushing old rows 10-51 down to 11-52
#
# Sheet "3_Industry_Analysis" (sheet3):
#   - rows 9/10/11 re-ranked now that "Foundry - Logic" has a 2nd member (UMC)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 1_Daily_Signals
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("1_Daily_Signals")

# Swap row 30 <-> row 31
$r30 = @($ws1.Cells.Item(30,1).Value(), $ws1.Cells.Item(30,2).Value(), $ws1.Cells.Item(30,3).Value(), $ws1.Cells.Item(30,4).Value(), $ws1.Cells.Item(30,5).Value(), $ws1.Cells.Item(30,6).Value(), $ws1.Cells.Item(30,7).Value(), $ws1.Cells.Item(30,8).Value())
$r31 = @($ws1.Cells.Item(31,1).Value(), $ws1.Cells.Item(31,2).Value(), $ws1.Cells.Item(31,3).Value(), $ws1.Cells.Item(31,4).Value(), $ws1.Cells.Item(31,5).Value(), $ws1.Cells.Item(31,6).Value(), $ws1.Cells.Item(31,7).Value(), $ws1.Cells.Item(31,8).Value())

$ws1.Cells.Item(30,1).Value = "'" + $r31[0]
$ws1.Cells.Item(30,2).Value = $r31[1]
$ws1.Cells.Item(30,3).Value = $r31[2]
$ws1.Cells.Item(30,4).Value = $r31[3]
$ws1.Cells.Item(30,5).Value = $r31[4]
$ws1.Cells.Item(30,6).Value = $r31[5]
$ws1.Cells.Item(30,7).Value = $r31[6]
$ws1.Cells.Item(30,8).Value = $r31[7]

$ws1.Cells.Item(31,1).Value = "'" + $r30[0]
$ws1.Cells.Item(31,2).Value = $r30[1]
$ws1.Cells.Item(31,3).Value = $r30[2]
$ws1.Cells.Item(31,4).Value = $r30[3]
$ws1.Cells.Item(31,5).Value = $r30[4]
$ws1.Cells.Item(31,6).Value = $r30[5]
$ws1.Cells.Item(31,7).Value = $r30[6]
$ws1.Cells.Item(31,8).Value = $r30[7]

# Swap row 36 <-> row 37
$r36 = @($ws1.Cells.Item(36,1).Value(), $ws1.Cells.Item(36,2).Value(), $ws1.Cells.Item(36,3).Value(), $ws1.Cells.Item(36,4).Value(), $ws1.Cells.Item(36,5).Value(), $ws1.Cells.Item(36,6).Value(), $ws1.Cells.Item(36,7).Value(), $ws1.Cells.Item(36,8).Value())
$r37 = @($ws1.Cells.Item(37,1).Value(), $ws1.Cells.Item(37,2).Value(), $ws1.Cells.Item(37,3).Value(), $ws1.Cells.Item(37,4).Value(), $ws1.Cells.Item(37,5).Value(), $ws1.Cells.Item(37,6).Value(), $ws1.Cells.Item(37,7).Value(), $ws1.Cells.Item(37,8).Value())

$ws1.Cells.Item(36,1).Value = "'" + $r37[0]
$ws1.Cells.Item(36,2).Value = $r37[1]
$ws1.Cells.Item(36,3).Value = $r37[2]
$ws1.Cells.Item(36,4).Value = $r37[3]
$ws1.Cells.Item(36,5).Value = $r37[4]
$ws1.Cells.Item(36,6).Value = $r37[5]
$ws1.Cells.Item(36,7).Value = $r37[6]
$ws1.Cells.Item(36,8).Value = $r37[7]

$ws1.Cells.Item(37,1).Value = "'" + $r36[0]
$ws1.Cells.Item(37,2).Value = $r36[1]
$ws1.Cells.Item(37,3).Value = $r36[2]
$ws1.Cells.Item(37,4).Value = $r36[3]
$ws1.Cells.Item(37,5).Value = $r36[4]
$ws1.Cells.Item(37,6).Value = $r36[5]
$ws1.Cells.Item(37,7).Value = $r36[6]
$ws1.Cells.Item(37,8).Value = $r36[7]

# New stock 2303 / UMC enters the daily-signals table -> insert a row at 47
# (old rows 47-51 shift down to 48-52)
$ws1.Rows.Item(47).Insert()
$ws1.Cells.Item(47,1).Value = "'2303"
$ws1.Cells.Item(47,2).Value = "聯電"
$ws1.Cells.Item(47,3).Value = "UMC"
$ws1.Cells.Item(47,4).Value = 62.8
$ws1.Cells.Item(47,5).Value = 0.64
$ws1.Cells.Item(47,6).Value = 34
$ws1.Cells.Item(47,7).Value = "Weak"
$ws1.Cells.Item(47,8).Value = 15.985

# ---------------------------------------------------------------------------
# Sheet 2: 2_21Day_Trend
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2_21Day_Trend")

# New stock 2303 / UMC ranks 10th by %_Tang_1_Thang -> insert a row at 10
# (old rows 10-51 shift down to 11-52)
$ws2.Rows.Item(10).Insert()
$ws2.Cells.Item(10,1).Value = "'2303"
$ws2.Cells.Item(10,2).Value = "聯電"
$ws2.Cells.Item(10,3).Value = "UMC"
$ws2.Cells.Item(10,4).Value = "Foundry - Logic"
$ws2.Cells.Item(10,5).Value = 16.3
$ws2.Cells.Item(10,6).Value = 0.42
$ws2.Cells.Item(10,7).Value = 15.985

# ---------------------------------------------------------------------------
# Sheet 3: 3_Industry_Analysis
# ---------------------------------------------------------------------------
# "Foundry - Logic" now has 2 members (TSMC + UMC); its aggregate stats move
# it above "Plastics" and "Power Supply", re-ranking rows 9-11.
$ws3 = $wb.Worksheets.Item("3_Industry_Analysis")

$ws3.Cells.Item(9,1).Value = "Foundry - Logic"
$ws3.Cells.Item(9,2).Value = 14.145
$ws3.Cells.Item(9,3).Value = 0.745
$ws3.Cells.Item(9,4).Value = 79.991
$ws3.Cells.Item(9,5).Value = 2

$ws3.Cells.Item(10,1).Value = "Plastics"
$ws3.Cells.Item(10,2).Value = 13.28
$ws3.Cells.Item(10,3).Value = 0.55
$ws3.Cells.Item(10,4).Value = 2.287
$ws3.Cells.Item(10,5).Value = 1

$ws3.Cells.Item(11,1).Value = "Power Supply"
$ws3.Cells.Item(11,2).Value = 13.18
$ws3.Cells.Item(11,3).Value = 0.98
$ws3.Cells.Item(11,4).Value = 20.343
$ws3.Cells.Item(11,5).Value = 2
